$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# Row 9
$ws.Range("D9").Value = "1399-05-12 (5)"
$ws.Range("E9").Value = "1400-04-15 (9)"
$ws.Range("F9").Value = "1401-04-05 (10)"
$ws.Range("G9").Value = "1402-02-27 (7)"
$ws.Range("H9").Value = "1402-02-27"

# Row 12
$ws.Range("D12").Value = 201327
$ws.Range("E12").Value = 109867
$ws.Range("F12").Value = 49576
$ws.Range("G12").Value = 409973
$ws.Range("H12").Value = 708620

# Row 13
$ws.Range("D13").Value = 119930
$ws.Range("E13").Value = 150930
$ws.Range("F13").Value = 342001
$ws.Range("G13").Value = 39251
$ws.Range("H13").Value = 11200

# Row 14
$ws.Range("D14").Value = 53085
$ws.Range("E14").Value = 144527
$ws.Range("F14").Value = 163778
$ws.Range("G14").Value = 1015661
$ws.Range("H14").Value = 2012247

# Row 15
$ws.Range("D15").Value = 84064
$ws.Range("E15").Value = 179362
$ws.Range("F15").Value = 1529149
$ws.Range("G15").Value = 2043642
$ws.Range("H15").Value = 4267680

# Row 16
$ws.Range("D16").Value = 64914
$ws.Range("E16").Value = 178769
$ws.Range("F16").Value = 459481
$ws.Range("G16").Value = 277476
$ws.Range("H16").Value = 244440

# Row 17
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0

# Row 18
$ws.Range("D18").Value = 523320
$ws.Range("E18").Value = 763455
$ws.Range("F18").Value = 2543985
$ws.Range("G18").Value = 3786003
$ws.Range("H18").Value = 7244187

# Row 19
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0

# Row 20
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 2034
$ws.Range("G20").Value = 2583
$ws.Range("H20").Value = 632

# Row 21
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0

# Row 22
$ws.Range("D22").Value = 673423
$ws.Range("E22").Value = 719533
$ws.Range("F22").Value = 723354
$ws.Range("G22").Value = 1016417
$ws.Range("H22").Value = 2522910

# Row 23
$ws.Range("D23").Value = 4591
$ws.Range("E23").Value = 4564
$ws.Range("F23").Value = 4737
$ws.Range("G23").Value = 4650
$ws.Range("H23").Value = 4564

# Row 24
$ws.Range("D24").Value = "-"
$ws.Range("E24").Value = "-"
$ws.Range("F24").Value = "-"
$ws.Range("G24").Value = "-"
$ws.Range("H24").Value = "-"

# Row 25
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0

# Row 26
$ws.Range("D26").Value = 678014
$ws.Range("E26").Value = 724097
$ws.Range("F26").Value = 730125
$ws.Range("G26").Value = 1023650
$ws.Range("H26").Value = 2528106

# Row 27
$ws.Range("D27").Value = 1201334
$ws.Range("E27").Value = 1487552
$ws.Range("F27").Value = 3274110
$ws.Range("G27").Value = 4809653
$ws.Range("H27").Value = 9772293

# Row 29
$ws.Range("D29").Value = 86732
$ws.Range("E29").Value = 320445
$ws.Range("F29").Value = 1800119
$ws.Range("G29").Value = 973521
$ws.Range("H29").Value = 1491363

# Row 30
$ws.Range("D30").Value = "-"
$ws.Range("E30").Value = "-"
$ws.Range("F30").Value = "-"
$ws.Range("G30").Value = "-"
$ws.Range("H30").Value = "-"

# Row 31
$ws.Range("D31").Value = 32181
$ws.Range("E31").Value = 56425
$ws.Range("F31").Value = 68615
$ws.Range("G31").Value = 136078
$ws.Range("H31").Value = 200353

# Row 32
$ws.Range("D32").Value = 34745
$ws.Range("E32").Value = 19655
$ws.Range("F32").Value = 27412
$ws.Range("G32").Value = 68882
$ws.Range("H32").Value = 86512

# Row 33
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 24624
$ws.Range("G33").Value = 71580
$ws.Range("H33").Value = 10306

# Row 34
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 1769990
$ws.Range("H34").Value = 5246746

# Row 35
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 0

# Row 36
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 0

# Row 37
$ws.Range("D37").Value = 153658
$ws.Range("E37").Value = 396525
$ws.Range("F37").Value = 1920770
$ws.Range("G37").Value = 3020051
$ws.Range("H37").Value = 7035280

# Row 38
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 0

# Row 39
$ws.Range("D39").Value = "-"
$ws.Range("E39").Value = "-"
$ws.Range("F39").Value = "-"
$ws.Range("G39").Value = "-"
$ws.Range("H39").Value = "-"

# Row 40
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 0

# Row 41
$ws.Range("D41").Value = 83068
$ws.Range("E41").Value = 103055
$ws.Range("F41").Value = 129316
$ws.Range("G41").Value = 133775
$ws.Range("H41").Value = 214532

# Row 42
$ws.Range("D42").Value = 83068
$ws.Range("E42").Value = 103055
$ws.Range("F42").Value = 129316
$ws.Range("G42").Value = 133775
$ws.Range("H42").Value = 214532

# Row 43
$ws.Range("D43").Value = 236726
$ws.Range("E43").Value = 499580
$ws.Range("F43").Value = 2050086
$ws.Range("G43").Value = 3153826
$ws.Range("H43").Value = 7249812

# Row 45
$ws.Range("D45").Value = 309166
$ws.Range("E45").Value = 309166
$ws.Range("F45").Value = 700000
$ws.Range("G45").Value = 700000
$ws.Range("H45").Value = 700000

# Row 46
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0

# Row 47
$ws.Range("D47").Value = 370781
$ws.Range("E47").Value = 390834
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0

# Row 48
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 0

# Row 49
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0

# Row 50
$ws.Range("D50").Value = 19083
$ws.Range("E50").Value = 30917
$ws.Range("F50").Value = 51719
$ws.Range("G50").Value = 70000
$ws.Range("H50").Value = 70000

# Row 51
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0

# Row 52
$ws.Range("D52").Value = "-"
$ws.Range("E52").Value = "-"
$ws.Range("F52").Value = "-"
$ws.Range("G52").Value = "-"
$ws.Range("H52").Value = "-"

# Row 53
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0

# Row 54
$ws.Range("D54").Value = "-"
$ws.Range("E54").Value = "-"
$ws.Range("F54").Value = "-"
$ws.Range("G54").Value = "-"
$ws.Range("H54").Value = "-"

# Row 55
$ws.Range("D55").Value = 0
$ws.Range("E55").Value = 0
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 0

# Row 56
$ws.Range("D56").Value = 265578
$ws.Range("E56").Value = 257055
$ws.Range("F56").Value = 472305
$ws.Range("G56").Value = 885827
$ws.Range("H56").Value = 1752481

# Row 57
$ws.Range("D57").Value = 964608
$ws.Range("E57").Value = 987972
$ws.Range("F57").Value = 1224024
$ws.Range("G57").Value = 1655827
$ws.Range("H57").Value = 2522481

# Row 58
$ws.Range("D58").Value = 1201334
$ws.Range("E58").Value = 1487552
$ws.Range("F58").Value = 3274110
$ws.Range("G58").Value = 4809653
$ws.Range("H58").Value = 9772293
